$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6, shifting existing rows 6-31 down to 7-32.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with its data.
$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(6, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(6, 4).Value = 44558
$ws.Cells.Item(6, 5).Value = 15
$ws.Cells.Item(6, 6).Value = "Fruta"
$ws.Cells.Item(6, 7).Value = 100103
$ws.Cells.Item(6, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(6, 9).Value = 100103004
$ws.Cells.Item(6, 10).Value = "Durazno"
$ws.Cells.Item(6, 11).Value = "Springcrest"
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 250
$ws.Cells.Item(6, 14).Value = 18000
$ws.Cells.Item(6, 15).Value = 20000
$ws.Cells.Item(6, 16).Value = 19000
$ws.Cells.Item(6, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(6, 18).Value = "Región Metropolitana"
$ws.Cells.Item(6, 19).Value = 1056
$ws.Cells.Item(6, 20).Value = 18
